# AutoCommit_13 июня 2024 г. 11:40:41_SibNout2023
# Update homework scores for a few students and move the frozen-pane
# viewport / active selection to reflect the newly edited rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (student #3): ДЗ_7 (J) raised from 2 to 5.
$ws.Range("J6").Value = 5

# Row 15 (student #12): ДЗ_5, ТК, ДЗ_6, ДЗ_7 (G:J) raised from 2 to 5.
$ws.Range("G15:J15").Value = 5

# Row 31 (student #28): ДЗ_5, ТК, ДЗ_6, ДЗ_7 (G:J) raised from 2 to 5.
$ws.Range("G31:J31").Value = 5

# Move the active selection to J31, matching where the edits were made
# (the frozen pane itself stays split at C/3 - only the active cell moves).
$ws.Activate()
[void]$ws.Range("J31").Select()
